# Update countries & provincias Spain
# Applies:
#   1. Swap country names "Finlandia" / "Serbia" (rows 43/44)
#   2. Swap country names "Principado de Andorra" / "Uzbekistan" (rows 84/85)
#   3. Update the "Datos actualizados..." timestamp (14:52 -> 15:22)
#   4. Update the updated case-count statistics for the affected rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1 & 2: swap the two pairs of country names ---------------------------
# Use temporary placeholder values so the swap doesn't collide mid-way.
$ws.Range("A43").Value = "__TMP_SWAP_1__"
$ws.Range("A44").Value = "Finlandia"
$ws.Range("A43").Value = "Serbia"

$ws.Range("A84").Value = "__TMP_SWAP_2__"
$ws.Range("A85").Value = "Principado de Andorra"
$ws.Range("A84").Value = "Uzbekistan"

# --- 3: update the timestamp text in A1 ------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 8 de Abril de 2020 a las 15:22"

# --- 4: update the numeric statistics ---------------------------------------
# Row 17 (Brasil)
$ws.Range("B17").Value = 14152
$ws.Range("C17").Value = 118
$ws.Range("E17").Value = 13326
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = 699

# Row 19 (Austria)
$ws.Range("B19").Value = 12852
$ws.Range("C19").Value = 213
$ws.Range("E19").Value = 8067

# Row 43 (now Serbia)
$ws.Range("B43").Value = 2666
$ws.Range("C43").Value = 219
$ws.Range("D43").Value = 118
$ws.Range("E43").Value = 2483
$ws.Range("F43").Value = 112
$ws.Range("G43").Value = 4
$ws.Range("H43").Value = 65

# Row 44 (now Finlandia)
$ws.Range("B44").Value = 2487
$ws.Range("C44").Value = 179
$ws.Range("D44").Value = 300
$ws.Range("E44").Value = 2147
$ws.Range("F44").Value = 82
$ws.Range("G44").Value = 6
$ws.Range("H44").Value = 40

# Row 79 (Republica de Macedonia)
$ws.Range("B79").Value = 617
$ws.Range("C79").Value = 18
$ws.Range("D79").Value = 35
$ws.Range("E79").Value = 553
$ws.Range("G79").Value = 3
$ws.Range("H79").Value = 29

# Row 84 (now Uzbekistan)
$ws.Range("C84").Value = 25
$ws.Range("D84").Value = 30
$ws.Range("E84").Value = 512
$ws.Range("F84").Value = 8
$ws.Range("G84").Value = 1
$ws.Range("H84").Value = 3

# Row 85 (now Principado de Andorra)
$ws.Range("B85").Value = 545
$ws.Range("C85").Value = 0
$ws.Range("D85").Value = 39
$ws.Range("E85").Value = 484
$ws.Range("F85").Value = 17
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 22
